# Updated symbol list on Thu Jan 26 17:29:54 UTC 2023 with GitHub Actions
#
# Refreshes the crypto price/volume table with new quotes. Every touched
# cell holds plain text (prices/percentages are stored as formatted text,
# not numbers) so we force the number format to Text ("@") before writing
# the value -- this stops Excel from re-interpreting strings like "304.31"
# or "1.02%" as numeric/percentage values -- and then restore the cell
# style to "Normal" so no stray style index is left behind on cells that
# originally had none.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, [string]$Ref, [string]$Val)
    $cell = $Sheet.Range($Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = "Normal"
}

# Row 2 - BNB
Set-TextCell $ws "D2" "304.31"
Set-TextCell $ws "E2" "1.02%"

# Row 3 - OKB
Set-TextCell $ws "D3" "35.85"
Set-TextCell $ws "E3" "1.53%"

# Row 4 - HuobiToken
Set-TextCell $ws "D4" "5.070"
Set-TextCell $ws "E4" "0.07%"

# Row 5 - Cronos
Set-TextCell $ws "E5" "1.72%"

# Row 6 - FTXToken
Set-TextCell $ws "D6" "1.934"
Set-TextCell $ws "E6" "2.50%"

# Row 7 - GateToken
Set-TextCell $ws "D7" "4.150"
Set-TextCell $ws "E7" "2.56%"

# Row 8 - KuCoinToken
Set-TextCell $ws "D8" "7.837"
Set-TextCell $ws "E8" "0.87%"

# Row 9 - MXToken
Set-TextCell $ws "D9" "0.9305"
Set-TextCell $ws "E9" "0.18%"

# Row 10 - LiechtensteinCryptoassetsExchange
Set-TextCell $ws "D10" "0.1292"
Set-TextCell $ws "E10" "-5.82%"

# Row 11 - WazirX
Set-TextCell $ws "D11" "0.1906"
Set-TextCell $ws "E11" "0.45%"

# Row 12 - MandalaExchangeToken
Set-TextCell $ws "D12" "0.09179"
Set-TextCell $ws "E12" "-0.30%"

# Row 13 - BitrueCoin
Set-TextCell $ws "D13" "0.03482"
Set-TextCell $ws "E13" "1.39%"

# Row 14 - BitMartToken
Set-TextCell $ws "D14" "0.09910"
Set-TextCell $ws "E14" "0.17%"

# Row 15 - BitForexToken
Set-TextCell $ws "D15" "0.001424"
Set-TextCell $ws "E15" "-0.29%"

# Row 16 - TigerCash
Set-TextCell $ws "D16" "0.006671"
Set-TextCell $ws "E16" "13.27%"

# Row 17 - LEO
Set-TextCell $ws "D17" "3.612"
Set-TextCell $ws "E17" "2.29%"

# Row 18 - BTSEToken
Set-TextCell $ws "D18" "3.039"
Set-TextCell $ws "E18" "3.88%"

# Row 19 - BitpandaEcosystemToken
Set-TextCell $ws "E19" "0.43%"

# Row 20 - was ProBitToken, now MCDex
Set-TextCell $ws "B20" "MCDex"
Set-TextCell $ws "C20" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell $ws "D20" "5.169"
Set-TextCell $ws "E20" "2.29%"

# Row 21 - was MCDex, now ProBitToken
Set-TextCell $ws "B21" "ProBitToken"
Set-TextCell $ws "C21" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell $ws "D21" "0.1303"
Set-TextCell $ws "E21" "0.17%"

# Row 22 - ZBToken
Set-TextCell $ws "D22" "0.2533"
Set-TextCell $ws "E22" "5.66%"

# Row 23 - CoinExToken
Set-TextCell $ws "D23" "0.04416"
Set-TextCell $ws "E23" "-1.81%"

# Row 24 - BitKan
Set-TextCell $ws "D24" "0.001236"
Set-TextCell $ws "E24" "1.77%"

# Row 25 - HotbitToken
Set-TextCell $ws "D25" "0.004695"
Set-TextCell $ws "E25" "-1.39%"

# Row 26 - NitroEx
Set-TextCell $ws "D26" "0.0001302"
Set-TextCell $ws "E26" "5.71%"

# Row 27 - UpBots
Set-TextCell $ws "D27" "0.0003133"
Set-TextCell $ws "E27" "4.29%"

# Row 39 - One
Set-TextCell $ws "D39" "0.01993"
Set-TextCell $ws "E39" "7.76%"

# Row 40 - IDEX
Set-TextCell $ws "D40" "0.05177"
Set-TextCell $ws "E40" "8.87%"

# Row 41 - KickToken
Set-TextCell $ws "D41" "0.007623"
Set-TextCell $ws "E41" "3.85%"

# Row 42 - Dexo
Set-TextCell $ws "D42" "0.01011"
Set-TextCell $ws "E42" "4.87%"

# Row 43 - BKEXToken
Set-TextCell $ws "D43" "0.1363"
Set-TextCell $ws "E43" "2.89%"

# Row 44 - CEJI
Set-TextCell $ws "D44" "0.002103"
Set-TextCell $ws "E44" "-0.46%"

# Row 45 - LocalTraders
Set-TextCell $ws "D45" "0.01073"
Set-TextCell $ws "E45" "-2.63%"

# Row 46 - CoinLion
Set-TextCell $ws "D46" "0.00006298"
Set-TextCell $ws "E46" "0.71%"

# Row 47 - Kangarootoken
Set-TextCell $ws "E47" "0.02%"

# Row 48 - BOLO
Set-TextCell $ws "D48" "65.22"
Set-TextCell $ws "E48" "0.84%"

# Row 49 - CoinbaseStockToken
Set-TextCell $ws "E49" "-3.45%"

# Row 50 - CryptobidCoin
Set-TextCell $ws "D50" "0.00002103"
Set-TextCell $ws "E50" "0.02%"

# Row 51 - SpecialPowerGold
Set-TextCell $ws "D51" "0.0002003"
Set-TextCell $ws "E51" "0.02%"
